# Auto-generated Excel COM-interop script applying the diff to before.xlsx
$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

$updates1 = @(
    @("C3", 1.270765083107522),
    @("D3", 0.9224612712860107),
    @("E3", 29.98),
    @("I3", 0.00485296642780304),
    @("J3", 59.96),
    @("C4", 0.8992907719883493),
    @("D4", 0.8475847244262695),
    @("E4", 33.23),
    @("I4", 0.06749633769989014),
    @("C5", 0.8083907055660962),
    @("D5", 0.6812993288040161),
    @("E5", 35.54),
    @("I5", 0.00437557235956192),
    @("J5", 66.45999999999999),
    @("C6", 0.7113278721406208),
    @("D6", 0.6888059973716736),
    @("E6", 36.63),
    @("I6", 0.07243952541351319),
    @("C7", 0.6330952457780761),
    @("D7", 0.6079980731010437),
    @("E7", 37.96),
    @("I7", 0.003739199161529541),
    @("J7", 71.08),
    @("C8", 0.5787000845118266),
    @("D8", 0.6008750796318054),
    @("E8", 38.23),
    @("I8", 0.0690688735961914),
    @("C9", 0.5223993065395975),
    @("D9", 0.5614283084869385),
    @("E9", 39.19),
    @("I9", 0.003415734624862671),
    @("J9", 73.26000000000001),
    @("C10", 0.4773743499100693),
    @("D10", 0.5118334293365479),
    @("E10", 39.29),
    @("I10", 0.07655083522796631),
    @("C11", 0.440242460103539),
    @("D11", 0.4911551177501678),
    @("E11", 39.42),
    @("I11", 0.00323596922159195),
    @("J11", 75.92),
    @("C12", 0.4073826791794319),
    @("D12", 0.519060492515564),
    @("E12", 40.13),
    @("I12", 0.07633972206115723),
    @("C13", 0.3596313433433936),
    @("D13", 0.5299490690231323),
    @("E13", 39.87),
    @("I13", 0.003082112908363342),
    @("J13", 76.45999999999999),
    @("C14", 0.3349607689836161),
    @("D14", 0.5293374061584473),
    @("E14", 40.04),
    @("I14", 0.0752610948562622),
    @("C15", 0.3008656945170426),
    @("D15", 0.5959169864654541),
    @("E15", 40.17),
    @("I15", 0.002885758626461029),
    @("J15", 78.38),
    @("C16", 0.2658336052322776),
    @("D16", 0.5468791723251343),
    @("E16", 40.73),
    @("I16", 0.07820269775390624),
    @("C17", 0.3532110272141976),
    @("D17", 0.5093264579772949),
    @("E17", 40.42),
    @("I17", 0.00283175522685051),
    @("J17", 78.58),
    @("C18", 0.3325512067573826),
    @("D18", 0.4972673356533051),
    @("E18", 40.72),
    @("I18", 0.08026079139709473),
    @("C19", 0.3196879969379766),
    @("D19", 0.5005775690078735),
    @("E19", 40.72),
    @("I19", 0.002830134373903274),
    @("J19", 78.84),
    @("C20", 0.3086711111349788),
    @("D20", 0.4958964288234711),
    @("E20", 40.75),
    @("I20", 0.07676160182952881),
    @("C21", 0.2967192929207794),
    @("D21", 0.5043570995330811),
    @("E21", 40.77),
    @("I21", 0.002720415234565735),
    @("J21", 80.26000000000001),
    @("C22", 0.3593374275095095),
    @("D22", 0.5436975359916687),
    @("E22", 40.54),
    @("I22", 0.08702284240722656),
    @("C23", 0.3407265352524393),
    @("D23", 0.5213198661804199),
    @("E23", 40.57),
    @("I23", 0.002815625762939453),
    @("J23", 79.73999999999999),
    @("C24", 0.3340221156434315),
    @("D24", 0.5743392109870911),
    @("E24", 40.5),
    @("I24", 0.08393898391723632),
    @("C25", 0.3298335441244327),
    @("D25", 0.5230526328086853),
    @("E25", 40.65),
    @("I25", 0.002781204640865326),
    @("J25", 80.08),
    @("C26", 0.3247077366927775),
    @("D26", 0.5035321712493896),
    @("E26", 40.59),
    @("I26", 0.07461999225616454),
    @("C27", 0.3721471850464984),
    @("D27", 0.5307528972625732),
    @("E27", 40.36),
    @("I27", 0.002891269946098328),
    @("J27", 80.34),
    @("C28", 0.3557094592389052),
    @("D28", 0.5453378558158875),
    @("E28", 40.46),
    @("I28", 0.09524719390869141),
    @("C29", 0.3493455309208816),
    @("D29", 0.5118885040283203),
    @("E29", 40.55),
    @("I29", 0.002815701484680176),
    @("J29", 81.45999999999999),
    @("C30", 0.345148349559404),
    @("D30", 0.5309150218963623),
    @("E30", 40.58),
    @("I30", 0.1016538452148437),
    @("C31", 0.3423345379713105),
    @("D31", 0.5333600044250488),
    @("E31", 40.5),
    @("I31", 0.002624308675527573),
    @("J31", 80.84),
    @("C32", 0.385869221474097),
    @("D32", 0.5219694972038269),
    @("E32", 40.02),
    @("I32", 0.0909037467956543),
    @("C33", 0.3710502609004819),
    @("D33", 0.4848226308822632),
    @("E33", 40.16),
    @("I33", 0.002589230251312256),
    @("J33", 81.44),
    @("C34", 0.3660709419386174),
    @("D34", 0.5567101240158081),
    @("E34", 40.22),
    @("I34", 0.09443283729553223),
    @("C35", 0.3626013264908054),
    @("D35", 0.5036542415618896),
    @("E35", 40.31),
    @("I35", 0.002615284967422485),
    @("J35", 81.44),
    @("C36", 0.3598312008671644),
    @("D36", 0.5082406401634216),
    @("E36", 40.38),
    @("I36", 0.09597752151489258),
    @("C37", 0.357073976983869),
    @("D37", 0.480644702911377),
    @("E37", 40.43),
    @("I37", 0.002610044485330582),
    @("J37", 81.5),
    @("C38", 0.3547871747637183),
    @("D38", 0.4923109710216522),
    @("E38", 40.44),
    @("I38", 0.0988216022491455),
    @("C39", 0.3532844435393326),
    @("D39", 0.4883528649806976),
    @("E39", 40.46),
    @("I39", 0.002656667786836624),
    @("J39", 81.54000000000001),
    @("A40", 1),
    @("B40", 38),
    @("C40", 0.3517466002121205),
    @("D40", 0.5189453363418579),
    @("E40", 40.48),
    @("I40", 0.09357086944580079),
    @("A41", 1),
    @("B41", 39),
    @("C41", 0.3502176706868458),
    @("D41", 0.5557814240455627),
    @("E41", 40.52),
    @("I41", 0.002609904646873474),
    @("J41", 81.08),
    @("A42", 1),
    @("B42", 40),
    @("C42", 0.3492598608741915),
    @("D42", 0.526050329208374),
    @("E42", 40.48),
    @("I42", 0.0878349796295166),
    @("I43", 0.002595285338163376),
    @("J43", 81.14),
    @("I44", 0.08899071083068848),
    @("I45", 0.002593597161769867),
    @("J45", 81),
    @("I46", 0.09076377182006835),
    @("I47", 0.002593524277210235),
    @("J47", 81.3),
    @("I48", 0.09262452125549317),
    @("I49", 0.002581756496429443),
    @("J49", 81.18000000000001),
    @("I50", 0.09330757446289062),
    @("I51", 0.002644850492477417),
    @("J51", 80.72),
    @("I52", 0.08167442455291749),
    @("I53", 0.002624920833110809),
    @("J53", 80.92),
    @("I54", 0.08468426246643067),
    @("I55", 0.002611336934566498),
    @("J55", 81.09999999999999),
    @("I56", 0.08634457855224609),
    @("I57", 0.002606158399581909),
    @("J57", 81.16),
    @("I58", 0.08763437843322754),
    @("I59", 0.002602447497844696),
    @("J59", 81),
    @("I60", 0.08843632011413574),
    @("I61", 0.002686895132064819),
    @("J61", 80.04000000000001),
    @("I62", 0.07816840209960937),
    @("I63", 0.002665933007001877),
    @("J63", 80.31999999999999),
    @("I64", 0.07939978275299073),
    @("I65", 0.002654549103975296),
    @("J65", 80.44),
    @("I66", 0.08044711017608643),
    @("I67", 0.002645803987979889),
    @("J67", 80.62),
    @("I68", 0.0813191785812378),
    @("I69", 0.002639099562168121),
    @("J69", 80.76000000000001),
    @("I70", 0.0821527910232544),
    @("I71", 0.002633226817846298),
    @("J71", 80.86),
    @("I72", 0.0829294719696045),
    @("I73", 0.002628669708967209),
    @("J73", 80.88),
    @("I74", 0.08363571052551269),
    @("I75", 0.002624688798189163),
    @("J75", 80.92),
    @("I76", 0.08421222915649414),
    @("G77", 1),
    @("H77", 38),
    @("I77", 0.002621688902378082),
    @("J77", 80.95999999999999),
    @("G78", 2),
    @("H78", 38),
    @("I78", 0.08475099678039551),
    @("J78", 0),
    @("G79", 1),
    @("H79", 39),
    @("I79", 0.002618488782644272),
    @("J79", 81.04000000000001),
    @("G80", 2),
    @("H80", 39),
    @("I80", 0.0851864917755127),
    @("J80", 0),
    @("G81", 1),
    @("H81", 40),
    @("I81", 0.002616441464424133),
    @("J81", 80.95999999999999),
    @("G82", 2),
    @("H82", 40),
    @("I82", 0.0856618408203125),
    @("J82", 0),
)

$updates2 = @(
    @("C3", 1.731323925049409),
    @("D3", 1.199891209602356),
    @("E3", 50.5),
    @("I3", 0.008701113414764405),
    @("J3", 29.78),
    @("C4", 1.056480479305205),
    @("D4", 1.005722999572754),
    @("E4", 60.05),
    @("I4", 0.004658307349681854),
    @("J4", 71.22),
    @("C5", 0.8920757427163746),
    @("D5", 0.8559004068374634),
    @("E5", 65.91),
    @("I5", 0.007340425610542298),
    @("J5", 44.18),
    @("C6", 0.7792850943363231),
    @("D6", 0.8242861032485962),
    @("E6", 68.53),
    @("I6", 0.003936619162559509),
    @("J6", 75.92),
    @("C7", 0.696297429178072),
    @("D7", 0.8060421347618103),
    @("E7", 69.25),
    @("I7", 0.006291057682037354),
    @("J7", 52.06),
    @("C8", 0.6190903245106988),
    @("D8", 0.7282993197441101),
    @("E8", 71.19),
    @("I8", 0.003266179448366165),
    @("J8", 79.76000000000001),
    @("C9", 0.5528593727427981),
    @("D9", 0.7849745154380798),
    @("E9", 71.40000000000001),
    @("I9", 0.005712466275691986),
    @("J9", 58.24),
    @("C10", 0.492773340124151),
    @("D10", 0.8573756217956543),
    @("E10", 69.81999999999999),
    @("I10", 0.003398492550849915),
    @("J10", 78.81999999999999),
    @("C11", 0.4243375367444495),
    @("D11", 0.7637041807174683),
    @("E11", 72.56999999999999),
    @("I11", 0.00580248544216156),
    @("J11", 59.24),
    @("C12", 0.3693596157042877),
    @("D12", 0.7921674251556396),
    @("E12", 71.86),
    @("I12", 0.00308498792052269),
    @("J12", 79.26000000000001),
    @("C13", 0.3053830882291431),
    @("D13", 0.8179534673690796),
    @("E13", 71.40000000000001),
    @("I13", 0.005602972781658173),
    @("J13", 60.32),
    @("C14", 0.4903627457502095),
    @("D14", 0.7236064672470093),
    @("E14", 71.84999999999999),
    @("I14", 0.002762129908800125),
    @("J14", 82.06),
    @("C15", 0.4607866257429123),
    @("D15", 0.7168929576873779),
    @("E15", 72.47),
    @("I15", 0.005895039045810699),
    @("J15", 60.52),
    @("C16", 0.443146137115748),
    @("D16", 0.7265926003456116),
    @("E16", 72.55),
    @("I16", 0.002779588043689728),
    @("J16", 82.28),
    @("C17", 0.4241101931294669),
    @("D17", 0.7443687319755554),
    @("E17", 72.5),
    @("I17", 0.005959517920017242),
    @("J17", 59.92),
    @("C18", 0.4094354672276456),
    @("D18", 0.7269152998924255),
    @("E18", 72.92),
    @("I18", 0.003255090802907944),
    @("J18", 79.72),
    @("C19", 0.391841823640077),
    @("D19", 0.7326924204826355),
    @("E19", 72.48999999999999),
    @("I19", 0.005835294270515442),
    @("J19", 63),
    @("C20", 0.3731644041836262),
    @("D20", 0.7307294607162476),
    @("E20", 72.39),
    @("I20", 0.002756731534004212),
    @("J20", 82.14),
    @("C21", 0.4298960707433846),
    @("D21", 0.717199444770813),
    @("E21", 72.59),
    @("I21", 0.006535295391082763),
    @("J21", 59.92),
    @("C22", 0.4242804021291111),
    @("D22", 0.723673403263092),
    @("E22", 72.5),
    @("I22", 0.002545238494873047),
    @("J22", 83.8),
    @("C23", 0.4199227324646452),
    @("D23", 0.7205944061279297),
    @("E23", 72.58),
    @("I23", 0.006682848465442657),
    @("J23", 61.4),
    @("C24", 0.4154623947713686),
    @("D24", 0.7242040634155273),
    @("E24", 72.61),
    @("I24", 0.002940166091918945),
    @("J24", 81.40000000000001),
    @("C25", 0.4109942846648071),
    @("D25", 0.7205256819725037),
    @("E25", 72.64),
    @("I25", 0.005749851787090301),
    @("J25", 60.8),
    @("C26", 0.4292265722609084),
    @("D26", 0.7179222702980042),
    @("E26", 72.59),
    @("I26", 0.002580700010061264),
    @("J26", 82.90000000000001),
    @("C27", 0.4260146538524524),
    @("D27", 0.7167750597000122),
    @("E27", 72.62),
    @("I27", 0.0054213014960289),
    @("J27", 63.02),
    @("C28", 0.4245100367976272),
    @("D28", 0.7229524850845337),
    @("E28", 72.65000000000001),
    @("I28", 0.002771724718809128),
    @("J28", 81.92),
    @("C29", 0.4233365358541841),
    @("D29", 0.721844494342804),
    @("E29", 72.59),
    @("I29", 0.005532943856716156),
    @("J29", 62.98),
    @("C30", 0.4218421439113824),
    @("D30", 0.7166802287101746),
    @("E30", 72.72),
    @("I30", 0.002677144342660904),
    @("J30", 82.12),
    @("C31", 0.4211825158285058),
    @("D31", 0.7290306687355042),
    @("E31", 72.64),
    @("I31", 0.005278353559970856),
    @("J31", 64.64),
    @("C32", 0.4197639523316984),
    @("D32", 0.7228803634643555),
    @("E32", 72.65000000000001),
    @("I32", 0.002942656654119492),
    @("J32", 80.36),
    @("C33", 0.4186395871574464),
    @("D33", 0.7260879278182983),
    @("E33", 72.63),
    @("I33", 0.005305182433128357),
    @("J33", 64.58),
    @("C34", 0.417834088044322),
    @("D34", 0.7277178764343262),
    @("E34", 72.61),
    @("I34", 0.002819725304841995),
    @("J34", 81.26000000000001),
    @("C35", 0.4165888240803843),
    @("D35", 0.7280434966087341),
    @("E35", 72.63),
    @("I35", 0.005476012766361236),
    @("J35", 63.48),
    @("C36", 0.4202570824519448),
    @("D36", 0.7194153070449829),
    @("E36", 72.67),
    @("I36", 0.002815248042345047),
    @("J36", 81.5),
    @("C37", 0.4199467030880244),
    @("D37", 0.7166546583175659),
    @("E37", 72.66),
    @("I37", 0.005891751098632813),
    @("J37", 62.32),
    @("C38", 0.4196536330425221),
    @("D38", 0.7187880873680115),
    @("E38", 72.64),
    @("I38", 0.002665816789865494),
    @("J38", 82.45999999999999),
    @("C39", 0.4194814990396085),
    @("D39", 0.7167866230010986),
    @("E39", 72.67),
    @("I39", 0.005458583652973175),
    @("J39", 63.08),
    @("C40", 0.4194517399629821),
    @("D40", 0.7178731560707092),
    @("E40", 72.68000000000001),
    @("I40", 0.002684184098243714),
    @("J40", 82.09999999999999),
    @("C41", 0.418865625301133),
    @("D41", 0.7199118733406067),
    @("E41", 72.66),
    @("I41", 0.005475252509117127),
    @("J41", 62.88),
    @("A42", 2),
    @("B42", 40),
    @("C42", 0.4187094795963038),
    @("D42", 0.7206975221633911),
    @("E42", 72.69999999999999),
    @("I42", 0.002702376937866211),
    @("J42", 82.12),
    @("I43", 0.005526572144031525),
    @("J43", 62.86),
    @("I44", 0.002680741119384765),
    @("J44", 82.3),
    @("I45", 0.005583512270450592),
    @("J45", 62.66),
    @("I46", 0.002658903688192367),
    @("J46", 82.56),
    @("I47", 0.005648680138587952),
    @("J47", 62.74),
    @("I48", 0.002631183755397796),
    @("J48", 82.54000000000001),
    @("I49", 0.005484263634681702),
    @("J49", 62.78),
    @("I50", 0.002659907591342926),
    @("J50", 82.40000000000001),
    @("I51", 0.005563467049598694),
    @("J51", 62.44),
    @("I52", 0.002631332474946976),
    @("J52", 82.8),
    @("I53", 0.005514841794967651),
    @("J53", 62.82),
    @("I54", 0.002658445787429809),
    @("J54", 82.48),
    @("I55", 0.005531346440315246),
    @("J55", 62.76),
    @("I56", 0.002655896151065827),
    @("J56", 82.42),
    @("I57", 0.005523120367527008),
    @("J57", 62.86),
    @("I58", 0.00265828327536583),
    @("J58", 82.58),
    @("I59", 0.00555011157989502),
    @("J59", 62.7),
    @("I60", 0.002654362088441849),
    @("J60", 82.58),
    @("I61", 0.005529772281646728),
    @("J61", 62.76),
    @("I62", 0.002665718126296997),
    @("J62", 82.54000000000001),
    @("I63", 0.005513874197006226),
    @("J63", 62.88),
    @("I64", 0.002674619770050049),
    @("J64", 82.38),
    @("I65", 0.005553813636302948),
    @("J65", 62.68),
    @("I66", 0.002653143012523651),
    @("J66", 82.54000000000001),
    @("I67", 0.005514228391647339),
    @("J67", 63.02),
    @("I68", 0.002681056839227676),
    @("J68", 82.23999999999999),
    @("I69", 0.005524485182762146),
    @("J69", 62.88),
    @("I70", 0.002658585298061371),
    @("J70", 82.45999999999999),
    @("I71", 0.005528382277488708),
    @("J71", 62.78),
    @("I72", 0.002658399724960327),
    @("J72", 82.54000000000001),
    @("I73", 0.00552231719493866),
    @("J73", 62.78),
    @("I74", 0.002662084347009659),
    @("J74", 82.5),
    @("I75", 0.005529314494132995),
    @("J75", 62.78),
    @("I76", 0.002658564758300781),
    @("J76", 82.56),
    @("I77", 0.005529915118217469),
    @("J77", 62.8),
    @("I78", 0.002658795362710953),
    @("J78", 82.56),
    @("I79", 0.005529847395420074),
    @("J79", 62.76),
    @("I80", 0.002659978419542313),
    @("J80", 82.56),
    @("G81", 1),
    @("H81", 40),
    @("I81", 0.005531894207000732),
    @("J81", 62.8),
    @("G82", 2),
    @("H82", 40),
    @("I82", 0.002660008710622787),
    @("J82", 82.59999999999999),
)

foreach ($u in $updates1) {
    $ws1.Range($u[0]).Value = $u[1]
}

foreach ($u in $updates2) {
    $ws2.Range($u[0]).Value = $u[1]
}
